$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename cluster label "MuSCs" -> "ECs" (shared string used by A2, A3, D2, D4)
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("D4").Value = "ECs"

# Updated numeric values (new TPM-derived stats)
$ws.Range("G2").Value = 0.1400336666666667
$ws.Range("H2").Value = 0.420101
$ws.Range("I2").Value = 0.001306326220890637
$ws.Range("J2").Value = 0.001306326220890637
$ws.Range("M2").Value = 0.1400336666666667
$ws.Range("N2").Value = 0.420101
$ws.Range("O2").Value = 0.001306326220890637
$ws.Range("P2").Value = 0.001306326220890637
$ws.Range("Q2").Value = 0.01960942780011111
$ws.Range("R2").Value = 0.176484850201
$ws.Range("S2").Value = 0.000001706488195386414
$ws.Range("T2").Value = 0.000001706488195386414

$ws.Range("G3").Value = 0.1400336666666667
$ws.Range("H3").Value = 0.420101
$ws.Range("I3").Value = 0.001306326220890637
$ws.Range("J3").Value = 0.001306326220890637
$ws.Range("M3").Value = 107.0565183333333
$ws.Range("N3").Value = 321.1695549999999
$ws.Range("O3").Value = 0.9986936737791093
$ws.Range("P3").Value = 0.9986936737791094
$ws.Range("Q3").Value = 14.99151680278389
$ws.Range("R3").Value = 134.923651225055
$ws.Range("S3").Value = 0.001304619732695251
$ws.Range("T3").Value = 0.001304619732695251

$ws.Range("G4").Value = 107.0565183333333
$ws.Range("H4").Value = 321.1695549999999
$ws.Range("I4").Value = 0.9986936737791093
$ws.Range("J4").Value = 0.9986936737791094
$ws.Range("M4").Value = 0.1400336666666667
$ws.Range("N4").Value = 0.420101
$ws.Range("O4").Value = 0.001306326220890637
$ws.Range("P4").Value = 0.001306326220890637
$ws.Range("Q4").Value = 14.99151680278389
$ws.Range("R4").Value = 134.923651225055
$ws.Range("S4").Value = 0.001304619732695251
$ws.Range("T4").Value = 0.001304619732695251

$ws.Range("G5").Value = 107.0565183333333
$ws.Range("H5").Value = 321.1695549999999
$ws.Range("I5").Value = 0.9986936737791093
$ws.Range("J5").Value = 0.9986936737791094
$ws.Range("M5").Value = 107.0565183333333
$ws.Range("N5").Value = 321.1695549999999
$ws.Range("O5").Value = 0.9986936737791093
$ws.Range("P5").Value = 0.9986936737791094
$ws.Range("Q5").Value = 11461.09811765533
$ws.Range("R5").Value = 103149.883058898
$ws.Range("S5").Value = 0.9973890540464141
$ws.Range("T5").Value = 0.9973890540464143
